$data = @(
    @('Natural Stone Inc', '1929 Dorsett, Houston, TX 77029'),
    @('Laminate Countertops', '3832 N Shepherd Dr B, Houston, TX 77018'),
    @('Tex-Star Marble and Granite', '4834 Ramus St, Houston, TX 77092'),
    @('BAYOU CITY GRANITE', '8307 Bauman Rd, Houston, TX 77022'),
    @('Superior Granite', '205 Meadow Lea Dr, Houston, TX 77022'),
    @('Stone Age Marble & Granite Inc', '4830 Ramus St, Houston, TX 77092'),
    @('AMC Granite & Marble', '3541 Creekmont Dr, Houston, TX 77091'),
    @('Aprile''s Marble & Granite', '6005 Milwee St # 714, Houston, TX 77092'),
    @('Francini Inc. Marble & Granite - Houston', '7108 Old Katy Rd suite 100, Houston, TX 77024'),
    @('CCH Granite', '3616 Lumberdale Rd, Houston, TX 77092'),
    @('Wilsonart', '552 Garden Oaks Blvd, Houston, TX 77018'),
    @('A.T. Custom Marble & Granite', '501 Gulf Bank Rd, Houston, TX 77037'),
    @('Terra Granite Slabs Center', '12642 Hempstead Rd, Houston, TX 77092'),
    @('Marquez Granite And Tile', '6826 Evans St, Houston, TX 77061'),
    @('Imperial Granite & Cabinets', '8700 Clay Rd #216, Houston, TX 77080'),
    @('Olympus Marble', '2141 Bingle Rd, Houston, TX 77055'),
    @('Granite Surfaces', '8030 Pinemont Dr, Houston, TX 77040'),
    @('Museum Flooring Granite Cabinet Factory', '4004 Gulf St, Houston, TX 77017'),
    @('Quality Countertops Inc', '1315 Illinois St, South Houston, TX 77587'),
    @('J.K. Granite & Cabinet', '6601 W Sam Houston Pkwy S, Houston, TX 77072'),
    @('Terraco Stone', '6411 Ashcroft Dr # A, Houston, TX 77081'),
    @('HTX Custom Countertops', '10580 Hammerly Blvd #1712, Houston, TX 77043'),
    @('MATERIAL Bespoke Stone + Tile', '1335 W Gray St Suite 100, Houston, TX 77019'),
    @('B and D Granite and Tile', '8010 Pinemont Dr, Houston, TX 77040'),
    @('Godstone Corp', '4722 Creekmont Dr, Houston, TX 77091'),
    @('International Granite & Marble (IGM)', '8017 Pinemont Dr # 300, Houston, TX 77040'),
    @('I Q Granite & Flooring Inc', '8017 Pinemont Dr # 300, Houston, TX 77040'),
    @('I Q Granite & Flooring Inc', '1306 Washington St, South Houston, TX 77587'),
    @('Toros Countertops & Gabinets', '1103 Bland St, Houston, TX 77091'),
    @('W R Watson Inc', '12902 Mula Ln, Stafford, TX 77477'),
    @('C M Custom Marble Design', '10129 O''Donnell Dr, Houston, TX 77076'),
    @('Vivaldi Lifestyles', '12614 Hempstead Rd, Houston, TX 77092'),
    @('Kitchen & Bath Concepts', '7026 Old Katy Rd #148, Houston, TX 77024'),
    @('Richmond Granite & Marble', '16403 W Bellfort Blvd, Sugar Land, TX 77498'),
    @('Space City Kitchen Remodeling Solutions', '2204a Louisiana St, Houston, TX 77002'),
    @('Roquemore Marble & Granite', '1291 N Post Oak Rd # 130, Houston, TX 77055'),
    @('The Home Depot', '1200 Home Depot Blvd, Sunset Valley, TX 78745'),
    @('Anna''s Granite Gallery', '16409 Farm to Market Rd 1325, Austin, TX 78728'),
    @('Capital of Texas Stoneworks', '9801 Beck Cir, Austin, TX 78758'),
    @('Lackstone Marbles & Granites', '2200 Denton Dr Suite-114, Austin, TX 78758'),
    @('CitiQuartz Austin LLC', '10609 Metric Blvd, Austin, TX 78758'),
    @('Langham Group', '10712 N Lamar Blvd c, Austin, TX 78753'),
    @('J & H Granite, Inc.', '237 Benelli Dr, Hutto, TX 78634'),
    @('J&L Granite Countertops', '135 Torchlight Dr, Bastrop, TX 78602'),
    @('PORCELANOSA Austin', '400 Barton Springs Rd, Austin, TX 78704'),
    @('Floor & Decor', '12901 N Interstate Hwy 35 #3, Austin, TX 78753'),
    @('The Home Depot', '1200 Barbara Jordan Blvd, Austin, TX 78723'),
    @('Stone Edge Countertops', '4502 Steffani Ln, Houston, TX 77041'),
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 479
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
